# Update the StimulationOrder sheet with the new table content and add a
# "blue" highlight fill (used for the second block of trials), mirroring
# the xlsx diff:
#   - A1 header changes from "overall trial" label position (unchanged text)
#   - Row values and wrapped shared strings change for rows 2-5
#   - Rows 4 & 5 (the second "block") get a new light-blue fill color
#   - Selection moves from D7 to C7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header) stays the same text --------------------------------
$ws.Range("A1").Value = "overall trial"
$ws.Range("B1").Value = "block"
$ws.Range("C1").Value = "trial"
$ws.Range("D1").Value = "channels"
$ws.Range("E1").Value = "electrodes"

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "[1, 2]"
$ws.Range("E2").Value = "[(1, 2), (3, 4)]"

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "[1, 8, 6, 5, 4, 7, 3]"
$ws.Range("E3").Value = "[(1, 2), (15, 16), (11, 12), (9, 10), (7, 8), (13, 14), (5, 6)]"

# --- Row 4 ---------------------------------------------------------------
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "[5, 7, 2, 6, 4, 1]"
$ws.Range("E4").Value = "[(9, 10), (13, 14), (3, 4), (11, 12), (7, 8), (1, 2)]"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "[2, 7, 4]"
$ws.Range("E5").Value = "[(3, 4), (13, 14), (7, 8)]"

# --- New fill colour for the second block (rows 4-5) ----------------------
# RGB(229, 229, 255) = 229 + 229*256 + 255*65536
$ws.Range("B4:E5").Interior.Color = 16770533

# --- Update the active selection to match the saved view ------------------
$ws.Range("C7").Select() | Out-Null
